$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.728.64"
$ws.Range("D3").Value = "1.895.29"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'248.46"
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.4943"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").Value = "'0.2970"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").Value = "'0.06826"
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("D10").Value = "1.897.77"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").Value = "'17.28"
$ws.Range("E11").Value = "  +3.15%  "
$ws.Range("D12").Value = "'92.51"
$ws.Range("E12").Value = "  +6.84%  "
$ws.Range("D13").Value = "'0.07263"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("E14").Value = "  +4.87%  "
$ws.Range("D15").Value = "'0.6812"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").Value = "30.706.95"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").Value = "'0.000007985"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("E18").Value = "  +4.02%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "2.140.91"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "'4.862"
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("D23").Value = "'192.97"
$ws.Range("E23").Value = "  +36.05%  "
$ws.Range("D24").Value = "'6.076"
$ws.Range("E24").Value = "  +6.17%  "
$ws.Range("D25").Value = "'9.432"
$ws.Range("E25").Value = "  +3.91%  "
$ws.Range("D26").Value = "'156.45"
$ws.Range("E26").Value = "  +4.34%  "
$ws.Range("D27").Value = "'19.29"
$ws.Range("E27").Value = "  +12.79%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").Value = "'1.404"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").Value = "'4.357"
$ws.Range("E30").Value = "  +3.77%  "
$ws.Range("D31").Value = "'0.09029"
$ws.Range("E31").Value = "  +2.93%  "
$ws.Range("D32").Value = "'4.045"
$ws.Range("E32").Value = "  +1.96%  "
$ws.Range("D33").Value = "'0.05212"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").Value = "'0.7470"
$ws.Range("E34").Value = "  +4.59%  "
$ws.Range("D35").Value = "'1.127"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").Value = "'2.735"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("D37").Value = "'0.01866"
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("D38").Value = "'2.689"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("D39").Value = "'2.170"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'0.9433"
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("D41").Value = "'0.4447"
$ws.Range("D42").Value = "'106.50"
$ws.Range("E42").Value = "  +4.18%  "
$ws.Range("D43").Value = "'5.787"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "'7.702"
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("E46").Value = "  +6.49%  "
$ws.Range("D47").Value = "'0.05864"
$ws.Range("E47").Value = "  +3.78%  "
$ws.Range("D48").Value = "'1.439"
$ws.Range("E48").Value = "  +7.37%  "
$ws.Range("D49").Value = "'8.735"
$ws.Range("E49").Value = "  +5.19%  "
$ws.Range("D50").Value = "'0.3966"
$ws.Range("E50").Value = "  +4.64%  "
$ws.Range("D51").Value = "'33.70"
$ws.Range("E51").Value = "  +3.76%  "
